$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: append a new sentence (as its own run) after "Linking erroneous
# fault codes to maintenance actions".
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Linking erroneous fault codes to maintenance actions" + [char]13) {
        $insertAt = $d.Range($p.Range.End - 1, $p.Range.End - 1)
        # Use a throwaway bookmark to force a fresh run boundary instead of
        # the new text being coalesced into the existing run.
        $d.Bookmarks.Add("zzzTempSplit", $insertAt)
        $insertAt.InsertAfter(" or finding the maintenance action that removes the correct fault code.")
        $d.Bookmarks.Item("zzzTempSplit").Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Change 2: after "3 codes may exist but one is an actual fault", add four
# new sub-bullets.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "3 codes may exist but one is an actual fault" + [char]13) {

        $p.Range.InsertParagraphAfter()
        $p1 = $d.Paragraphs.Item($i + 1)
        $p1.Range.Text = "Filter by corrosion/bare metal/ corrosion Preventative treatment"

        $p1.Range.InsertParagraphAfter()
        $p2 = $d.Paragraphs.Item($i + 2)
        $p2.Range.ListFormat.ListIndent()
        $p2.Range.Text = "See what maintenance action codes (action code taken) "

        $p2.Range.InsertParagraphAfter()
        $p3 = $d.Paragraphs.Item($i + 3)
        $p3.Range.ListFormat.ListOutdent()
        $p3.Range.Text = "See what MSP codes exist between (received and completion date) and which do not shortly after (hyper-parameter we will tune). "

        $p3.Range.InsertParagraphAfter()
        $p4 = $d.Paragraphs.Item($i + 4)
        $p4.Range.ListFormat.ListIndent()
        $p4.Range.Text = "Potentially list first 10 that drop off"

        break
    }
}

# ---------------------------------------------------------------------------
# Change 3: split "Clustering algorithms" into "Clustering algo" + "rithms",
# relocating the "_GoBack" bookmark to sit between the two new runs.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Clustering algorithms" + [char]13) {
        $splitAt = $p.Range.Start + ("Clustering algo").Length
        $splitRange = $d.Range($splitAt, $splitAt)
        $d.Bookmarks.Add("_GoBack", $splitRange)
        break
    }
}
